$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.740.22"
$ws.Range("E2").Value = "  +6.16%  "

# Row 3
$ws.Range("D3").Value = "2.055.05"
$ws.Range("E3").Value = "  +3.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'252.86"
$ws.Range("E5").Value = "  +3.58%  "

# Row 6
$ws.Range("D6").Value = "'0.649"
$ws.Range("E6").Value = "  +2.59%  "

# Row 7
$ws.Range("D7").Value = "'65.35"
$ws.Range("E7").Value = "  +14.09%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "'60.75"
$ws.Range("E9").Value = "  +2.35%  "

# Row 10
$ws.Range("E10").Value = "  +5.76%  "

# Row 11
$ws.Range("D11").Value = "'0.0765"
$ws.Range("E11").Value = "  +4.90%  "

# Row 12
$ws.Range("D12").Value = "'0.105"
$ws.Range("E12").Value = "  +1.65%  "

# Row 13
$ws.Range("D13").Value = "'0.919"
$ws.Range("E13").Value = "  -1.18%  "

# Row 14
$ws.Range("D14").Value = "'14.98"
$ws.Range("E14").Value = "  +5.44%  "

# Row 15
$ws.Range("D15").Value = "2.357.39"
$ws.Range("E15").Value = "  +3.57%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'5.62"
$ws.Range("E16").Value = "  +7.48%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'20.70"
$ws.Range("E17").Value = "  +20.64%  "

# Row 18
$ws.Range("D18").Value = "2.049.33"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19
$ws.Range("D19").Value = "37.672.83"
$ws.Range("E19").Value = "  +6.04%  "

# Row 20
$ws.Range("D20").Value = "'73.96"
$ws.Range("E20").Value = "  +4.37%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  +4.82%  "

# Row 22
$ws.Range("D22").Value = "'5.43"
$ws.Range("E22").Value = "  +5.63%  "

# Row 23
$ws.Range("D23").Value = "'241.18"
$ws.Range("E23").Value = "  +3.37%  "

# Row 24
$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +8.41%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +3.75%  "

# Row 27
$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = "  +5.88%  "

# Row 28
$ws.Range("D28").Value = "'162.18"
$ws.Range("E28").Value = "  -0.69%  "

# Row 29
$ws.Range("D29").Value = "'20.00"
$ws.Range("E29").Value = "  +3.39%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'5.36"
$ws.Range("E30").Value = "  +11.47%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.117"
$ws.Range("E31").Value = "  +30.56%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = "  +9.42%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  +2.78%  "

# Row 34
$ws.Range("D34").Value = "'4.75"
$ws.Range("E34").Value = "  +12.52%  "

# Row 35
$ws.Range("D35").Value = "'0.0625"
$ws.Range("E35").Value = "  +5.42%  "

# Row 36
$ws.Range("E36").Value = "  +3.20%  "

# Row 37
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "  +3.53%  "

# Row 39
$ws.Range("D39").Value = "'6.14"
$ws.Range("E39").Value = "  +24.64%  "

# Row 40
$ws.Range("E40").Value = "  +17.84%  "

# Row 41
$ws.Range("D41").Value = "'2.88"
$ws.Range("E41").Value = "  +26.44%  "

# Row 42
$ws.Range("D42").Value = "'1.24"
$ws.Range("E42").Value = "  +4.62%  "

# Row 43
$ws.Range("D43").Value = "'0.0221"
$ws.Range("E43").Value = "  +4.67%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'17.29"
$ws.Range("E44").Value = "  +9.55%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.15"
$ws.Range("E45").Value = "  +5.82%  "

# Row 46
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.92"
$ws.Range("E46").Value = "  +2.79%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'95.74"
$ws.Range("E47").Value = "  +5.12%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.95"
$ws.Range("E48").Value = "  +6.77%  "

# Row 49
$ws.Range("D49").Value = "1.405.82"
$ws.Range("E49").Value = "  +0.79%  "

# Row 50
$ws.Range("E50").Value = "  +2.03%  "

# Row 51
$ws.Range("D51").Value = "'47.07"
$ws.Range("E51").Value = "  +4.32%  "
